$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.08464515209198
$ws.Range("B1").Value = 2.281284809112549
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.512775182723999
$ws.Range("E1").Value = 0.9896021485328674
